$d = $word.ActiveDocument

$old = "Ви сте учесници глобалног посматрачког пројекта, који има за циљ да одреди колико је светлосно загађене у средини у којој живите. Посматрајући звезде унутар сазвежђа Персеус и упоређујући их са приложеним звезданим картама, посматрачи широм света могу на практичном примеру да увиде колико је светлосно загађење у њиховој средини. Кроз учешће у овом пројекту, допринећете целовитијем сагледавању глобалног проблема."
$new = "Ви сте учесници глобалног посматрачког пројекта, који има за циљ да одреди колико је светлосно загађене у средини у којој живите. Посматрајући звезде унутар  Херкулово сазвежђе и упоређујући их са приложеним звезданим картама, посматрачи широм света могу на практичном примеру да увиде колико је светлосно загађење у њиховој средини. Кроз учешће у овом пројекту, допринећете целовитијем сагледавању глобалног проблема."

foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    if ($pr.Text.TrimEnd([char]13) -eq $old) {
        $start = $pr.Start
        $r = $d.Range($start, $pr.End)
        $r.Text = ""
        $r2 = $d.Range($start, $start)
        $r2.InsertAfter($new)
        break
    }
}
